$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-03-09 Saturday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-03-10 Sunday", 2) | Out-Null
$d.Content.Find.Execute("240÷5=48, 0", $true, $false, $false, $false, $false, $true, 1, $false, "705÷7=100, 5", 2) | Out-Null
$d.Content.Find.Execute("532÷8=66, 4", $true, $false, $false, $false, $false, $true, 1, $false, "814÷5=162, 4", 2) | Out-Null
$d.Content.Find.Execute("826÷3=275, 1", $true, $false, $false, $false, $false, $true, 1, $false, "978÷5=195, 3", 2) | Out-Null
$d.Content.Find.Execute("722÷4=180, 2", $true, $false, $false, $false, $false, $true, 1, $false, "834÷4=208, 2", 2) | Out-Null
$d.Content.Find.Execute("546÷3=182, 0", $true, $false, $false, $false, $false, $true, 1, $false, "851÷8=106, 3", 2) | Out-Null
$d.Content.Find.Execute("522÷7=74, 4", $true, $false, $false, $false, $false, $true, 1, $false, "126÷2=63, 0", 2) | Out-Null
$d.Content.Find.Execute("497÷7=71, 0", $true, $false, $false, $false, $false, $true, 1, $false, "906÷4=226, 2", 2) | Out-Null
$d.Content.Find.Execute("156÷5=31, 1", $true, $false, $false, $false, $false, $true, 1, $false, "918÷7=131, 1", 2) | Out-Null
$d.Content.Find.Execute("659÷5=131, 4", $true, $false, $false, $false, $false, $true, 1, $false, "794÷4=198, 2", 2) | Out-Null
$d.Content.Find.Execute("191÷4=47, 3", $true, $false, $false, $false, $false, $true, 1, $false, "630÷2=315, 0", 2) | Out-Null
$d.Content.Find.Execute("196÷9=21, 7", $true, $false, $false, $false, $false, $true, 1, $false, "428÷8=53, 4", 2) | Out-Null
$d.Content.Find.Execute("363÷2=181, 1", $true, $false, $false, $false, $false, $true, 1, $false, "173÷8=21, 5", 2) | Out-Null
$d.Content.Find.Execute("357÷9=39, 6", $true, $false, $false, $false, $false, $true, 1, $false, "394÷4=98, 2", 2) | Out-Null
$d.Content.Find.Execute("377÷7=53, 6", $true, $false, $false, $false, $false, $true, 1, $false, "135÷6=22, 3", 2) | Out-Null
$d.Content.Find.Execute("962÷5=192, 2", $true, $false, $false, $false, $false, $true, 1, $false, "629÷9=69, 8", 2) | Out-Null
$d.Content.Find.Execute("709÷9=78, 7", $true, $false, $false, $false, $false, $true, 1, $false, "101÷4=25, 1", 2) | Out-Null
$d.Content.Find.Execute("644÷8=80, 4", $true, $false, $false, $false, $false, $true, 1, $false, "411÷3=137, 0", 2) | Out-Null
$d.Content.Find.Execute("312÷8=39, 0", $true, $false, $false, $false, $false, $true, 1, $false, "843÷2=421, 1", 2) | Out-Null
$d.Content.Find.Execute("525÷3=175, 0", $true, $false, $false, $false, $false, $true, 1, $false, "350÷2=175, 0", 2) | Out-Null
$d.Content.Find.Execute("770÷3=256, 2", $true, $false, $false, $false, $false, $true, 1, $false, "845÷4=211, 1", 2) | Out-Null
$d.Content.Find.Execute("766÷8=95, 6", $true, $false, $false, $false, $false, $true, 1, $false, "376÷7=53, 5", 2) | Out-Null
$d.Content.Find.Execute("374÷7=53, 3", $true, $false, $false, $false, $false, $true, 1, $false, "880÷5=176, 0", 2) | Out-Null
$d.Content.Find.Execute("592÷7=84, 4", $true, $false, $false, $false, $false, $true, 1, $false, "840÷2=420, 0", 2) | Out-Null
$d.Content.Find.Execute("492÷9=54, 6", $true, $false, $false, $false, $false, $true, 1, $false, "658÷5=131, 3", 2) | Out-Null
$d.Content.Find.Execute("321÷5=64, 1", $true, $false, $false, $false, $false, $true, 1, $false, "826÷6=137, 4", 2) | Out-Null
